$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-08-07 Wednesday" "2024-08-08 Thursday"

Replace-Text "800×9=7200" "811×2=1622"
Replace-Text "197×5=985" "749×4=2996"
Replace-Text "763×2=1526" "146×9=1314"
Replace-Text "193×9=1737" "832×7=5824"
Replace-Text "831×7=5817" "418×5=2090"

Replace-Text "867×2=1734" "253×7=1771"
Replace-Text "954×3=2862" "353×5=1765"
Replace-Text "690×3=2070" "788×3=2364"
Replace-Text "589×7=4123" "429×7=3003"
Replace-Text "925×5=4625" "616×2=1232"

Replace-Text "266×2=532" "602×9=5418"
Replace-Text "197×9=1773" "168×6=1008"
Replace-Text "114×9=1026" "512×9=4608"
Replace-Text "299×8=2392" "664×5=3320"
Replace-Text "711×2=1422" "385×4=1540"

Replace-Text "220×2=440" "673×9=6057"
Replace-Text "601×6=3606" "139×9=1251"
Replace-Text "172×7=1204" "444×5=2220"
Replace-Text "448×9=4032" "882×5=4410"
Replace-Text "717×3=2151" "498×7=3486"

Replace-Text "134×9=1206" "532×8=4256"
Replace-Text "864×7=6048" "682×7=4774"
Replace-Text "388×9=3492" "824×7=5768"
Replace-Text "166×9=1494" "389×4=1556"
Replace-Text "243×9=2187" "193×8=1544"
